# Updated cryptos list on Mon Apr 24 14:44:22 UTC 2023 with GitHub Actions
#
# The sheet stores "Price" (col D) and "Volume(1h)" (col E) as plain text,
# even though many of the price strings (e.g. "336.90", "1.012") look like
# numbers. Writing straight to .Value would make Excel reinterpret them as
# numeric values (dropping trailing zeros, stripping the padding spaces /
# percent sign handling, etc.) and - via NumberFormat - would also leave a
# stray text-format style on the cell. To keep both the literal text and
# the original (unstyled) cell formatting intact, force the cell to Text
# format, assign the string, then restore the cell style to Normal (which
# also resets NumberFormat back to General without touching the value).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, [string]$val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# Row 51 coin/link change: EOS -> PancakeSwap
Set-TextValue $ws.Cells.Item(51, 2) "PancakeSwap"
Set-TextValue $ws.Cells.Item(51, 3) "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"

# Updated Price (column D) and Volume(1h) (column E) values for rows 2-51.
$updates = @(
    @{Row=2;  D="27.694.32";  E="  -0.54%  "},
    @{Row=3;  D="1.876.65";   E="  -0.55%  "},
    @{Row=4;  D="1.015";      E="  +0.70%  "},
    @{Row=5;  D="336.13";     E="  +0.59%  "},
    @{Row=6;  D="1.012";      E="  +0.49%  "},
    @{Row=7;  D="0.4649";     E="  -1.89%  "},
    @{Row=8;  D="0.3942";     E="  +0.31%  "},
    @{Row=9;  D="45.89";      E="  -3.84%  "},
    @{Row=10; D="0.07991";    E="  -1.02%  "},
    @{Row=11; D="1.008";      E="  -1.78%  "},
    @{Row=12; D="21.90";      E="  -1.10%  "},
    @{Row=13; D="1.886.93";   E="  +0.07%  "},
    @{Row=14; D="5.980";      E="  -0.13%  "},
    @{Row=15; D="7.259";      E="  +1.54%  "},
    @{Row=16; D="1.016";      E="  +0.44%  "},
    @{Row=17; D="88.98";      E="  +1.82%  "},
    @{Row=18; D="0.06727";    E="  +0.05%  "},
    @{Row=19; D="0.00001046"; E="  -0.49%  "},
    @{Row=20; D="17.27";      E="  -0.33%  "},
    @{Row=21; D="1.009";      E="  +0.18%  "},
    @{Row=22; D="27.744.19";  E="  -0.49%  "},
    @{Row=23; D="5.474";      E="  -0.85%  "},
    @{Row=24; D="10.97";      E="  -0.13%  "},
    @{Row=25; D="2.304";      E="  -1.20%  "},
    @{Row=26; D="2.106.50";   E="  -0.20%  "},
    @{Row=27; D="158.87";     E="  -0.23%  "},
    @{Row=28; D="19.75";      E="  -1.66%  "},
    @{Row=29; D="2.147";      E="  +1.96%  "},
    @{Row=30; D="5.473";      E="  -1.69%  "},
    @{Row=31; D="121.54";     E="  -0.43%  "},
    @{Row=32; D="0.9803";     E="  +0.25%  "},
    @{Row=33; D="0.09438";    E="  -0.67%  "},
    @{Row=34; D="3.630";      E="  -0.07%  "},
    @{Row=35; D="5.316";      E="  -0.82%  "},
    @{Row=36; D="1.346";      E="  -7.17%  "},
    @{Row=37; D="0.06061";    E="  -1.74%  "},
    @{Row=38; D="0.02237";    E="  -1.55%  "},
    @{Row=39; D="8.342";      E="  +3.32%  "},
    @{Row=40; D="1.199";      E="  -1.77%  "},
    @{Row=41; D="1.011";      E="  +0.41%  "},
    @{Row=42; D="0.5969";     E="  -0.64%  "},
    @{Row=43; D="0.1890";     E="  -0.47%  "},
    @{Row=44; D="10.36";      E="  +0.32%  "},
    @{Row=45; D="1.247";      E="  -0.98%  "},
    @{Row=46; D="0.5641";     E="  -1.17%  "},
    @{Row=47; D="12.28";      E="  -0.03%  "},
    @{Row=48; D="1.934";      E="  -0.35%  "},
    @{Row=49; D="0.06776";    E="  -2.00%  "},
    @{Row=50; D="111.79";     E="  -1.56%  "},
    @{Row=51; D="3.017";      E="  -11.33%  "}
)

foreach ($u in $updates) {
    Set-TextValue $ws.Cells.Item($u.Row, 4) $u.D
    Set-TextValue $ws.Cells.Item($u.Row, 5) $u.E
}
